# Generate Report for Handoff
#
# c44e843b-c8b2-4011-a08e-70b02a849b95.md moved from "In Translation" to
# "Ready for handoff", and the status-check timestamp was refreshed for the
# three rows that have no completed handoff yet (rows 6, 9, 10 - the ones
# missing a "Latest Target File"/"Latest Handback File") on every sheet.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D6").Value  = "2016-18-11 20:18:05"
$wsOverview.Range("B9").Value  = "Ready for handoff"
$wsOverview.Range("C9").Value  = "Ready for handoff"
$wsOverview.Range("D9").Value  = "2016-18-11 20:18:05"
$wsOverview.Range("D10").Value = "2016-18-11 20:18:05"

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E6").Value  = "2016-03-11 20:18:01"
$wsZhCn.Range("C9").Value  = "Ready for handoff"
$wsZhCn.Range("E9").Value  = "2016-03-11 20:18:01"
$wsZhCn.Range("E10").Value = "2016-03-11 20:18:01"

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E6").Value  = "2016-03-11 20:18:05"
$wsDeDe.Range("C9").Value  = "Ready for handoff"
$wsDeDe.Range("E9").Value  = "2016-03-11 20:18:05"
$wsDeDe.Range("E10").Value = "2016-03-11 20:18:05"
